$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = "Adding Conflict to git"
$ws.Range("G9").Value = "2nd Conflict"

$ws.Range("G11").Select()
